$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell as exact text (avoids Excel auto-converting numeric-looking strings)
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 12 <-> Row 13: Cardano/Toncoin order swapped
$ws.Range("B12").Value = "Toncoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D12") "5.13"
$ws.Range("E12").Value = "  -2.03%  "

$ws.Range("B13").Value = "Cardano"
$ws.Range("C13").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue $ws.Range("D13") "0.349"
$ws.Range("E13").Value = "  -4.37%  "

# Per-row price/volume updates
Set-TextValue $ws.Range("D2") "66.610.70"
$ws.Range("E2").Value = "  -1.53%  "
Set-TextValue $ws.Range("D3") "2.516.89"
$ws.Range("E3").Value = "  -4.49%  "
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.09%  "
Set-TextValue $ws.Range("D5") "584.17"
$ws.Range("E5").Value = "  -1.76%  "
Set-TextValue $ws.Range("D6") "171.58"
$ws.Range("E6").Value = "  +1.85%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -1.57%  "
Set-TextValue $ws.Range("D9") "2.514.74"
$ws.Range("E9").Value = "  -4.54%  "
$ws.Range("E10").Value = "  -1.28%  "
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("E14").Value = "  -3.57%  "
Set-TextValue $ws.Range("D15") "2.974.10"
$ws.Range("E15").Value = "  -4.60%  "
$ws.Range("E16").Value = "  -3.41%  "
Set-TextValue $ws.Range("D17") "66.448.80"
$ws.Range("E17").Value = "  -1.70%  "
Set-TextValue $ws.Range("D18") "2.514.67"
$ws.Range("E18").Value = "  -3.99%  "
Set-TextValue $ws.Range("D19") "7.84"
$ws.Range("E19").Value = "  -2.74%  "
Set-TextValue $ws.Range("D20") "11.29"
$ws.Range("E20").Value = "  -6.25%  "
Set-TextValue $ws.Range("D21") "347.73"
$ws.Range("E21").Value = "  -2.87%  "
$ws.Range("E22").Value = "  -3.03%  "
$ws.Range("E23").Value = "  -1.26%  "
Set-TextValue $ws.Range("D24") "1.98"
$ws.Range("E24").Value = "  +1.70%  "
$ws.Range("E25").Value = "  -0.15%  "
Set-TextValue $ws.Range("D26") "70.28"
$ws.Range("E26").Value = "  +0.33%  "
Set-TextValue $ws.Range("D27") "9.97"
$ws.Range("E27").Value = "  -3.95%  "
Set-TextValue $ws.Range("D28") "1.00"
$ws.Range("E28").Value = "  -0.58%  "
Set-TextValue $ws.Range("D29") "2.631.96"
$ws.Range("E29").Value = "  -4.98%  "
Set-TextValue $ws.Range("D30") "0.0₃0977"
$ws.Range("E30").Value = "  -3.37%  "
Set-TextValue $ws.Range("D31") "525.76"
$ws.Range("E31").Value = "  -4.40%  "
Set-TextValue $ws.Range("D32") "8.10"
$ws.Range("E32").Value = "  +1.67%  "
$ws.Range("E33").Value = "  -3.07%  "
$ws.Range("E34").Value = "  -3.23%  "
$ws.Range("E35").Value = "  -4.91%  "
Set-TextValue $ws.Range("D36") "0.999"
$ws.Range("E36").Value = "  -0.11%  "
Set-TextValue $ws.Range("D37") "1.47"
$ws.Range("E37").Value = "  -2.99%  "
Set-TextValue $ws.Range("D38") "156.98"
$ws.Range("E38").Value = "  -0.85%  "
Set-TextValue $ws.Range("D39") "18.64"
$ws.Range("E39").Value = "  -2.36%  "
$ws.Range("E40").Value = "  +0.58%  "
$ws.Range("E41").Value = "  -3.25%  "
$ws.Range("E42").Value = "  -1.37%  "
$ws.Range("E43").Value = "  -2.50%  "
$ws.Range("E44").Value = "  -0.09%  "
Set-TextValue $ws.Range("D45") "2.50"
$ws.Range("E45").Value = "  +2.58%  "
$ws.Range("E46").Value = "  -1.60%  "
Set-TextValue $ws.Range("D47") "149.20"
$ws.Range("E47").Value = "  -2.56%  "
Set-TextValue $ws.Range("D48") "0.559"
$ws.Range("E48").Value = "  -4.04%  "
$ws.Range("E49").Value = "  -3.35%  "
$ws.Range("E50").Value = "  +0.39%  "
$ws.Range("E51").Value = "  -11.11%  "
